# The "Population" figures in column H (rows 66-119 of the
# "Quarterly Data Levels" sheet) were off by one row: every value from H66
# down to H118 actually belonged to the row below it, and the final row's
# figure (H118/H119 = "5" / "380") was really a single value, 5380, that had
# been mistakenly split across two cells.
#
# Fix: shift H66:H118 down by one row (H67 <- old H66, H68 <- old H67, ...,
# H118 <- old H117), clear out H66 (nothing should shift into it), and set
# H119 to the corrected value 5380. Every other cell/column is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly Data Levels")

$firstRow = 66
$lastRow = 118
$col = 8  # column H

# Snapshot the original values for H66:H118 before we overwrite anything,
# using Value() (a method call) since Value as a bare property returns a
# descriptor rather than the actual cell contents in this host.
$origValues = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $origValues[$r] = $ws.Cells.Item($r, $col).Value()
}

# Push every value down into the row below it: H67 = old H66, H68 = old H67,
# ..., H118 = old H117.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $ws.Cells.Item($r, $col).Value = $origValues[$r - 1]
}

# The old H66 value has now moved to H67, so H66 itself becomes blank.
$ws.Cells.Item($firstRow, $col).ClearContents()

# H119 previously held "380" (paired with H118's stray "5"); replace it with
# the corrected combined figure.
$ws.Cells.Item($lastRow + 1, $col).Value = 5380
